# Lowercase the "exercises/Exercise_*.html" links to "exercises/exercise_*.html"
# throughout the "In-Class Exercise" column (column D) of the class schedule.
#
# Cells in this sheet that hold these backtick-led RST links carry a
# "quote prefix" style (so Excel doesn't try to parse the leading backtick
# as a formula). A plain `.Value` assignment would normally cause the
# style to lose that quote-prefix flag, so we re-assert it by leading the
# new value with a literal apostrophe -- exactly like a user typing
# '`Ex ... in the Excel UI -- which keeps the original cell style (s=)
# intact and only changes the shared-string text/index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2
    if ($v -and $v -match "Exercise_") {
        $newValue = $v -replace "Exercise_", "exercise_"
        $cell.Value = "'" + $newValue
    }
}

# Move the active selection from C12 to C11, matching the saved workbook state.
$ws.Range("C11").Select()
